# Applies the "equiblissExt custom scheduler" edit:
#  - flips the Max-Slowdown (row 8) / Weighted-Speedup (row 9) formulas to the
#    new orientation (ratio inverted) to support the new custom scheduler data
#  - repoints P:S (the "4-wide" summary columns) to simply mirror H:K instead
#    of recomputing independently
#  - inserts a new row 11 ("Fairness" = WghtSpeedup / MaxSlowdown)
#  - adds a new "IPC alone vs IPC shared" scratch calculation block (rows 63-72)
#  - updates the 7 charts whose series ranges shift down by one row
#  - updates the saved selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 8 ("Max Slow"): was MAX($C$x/col x) now MAX(col x/$C$x)
# ---------------------------------------------------------------------------
$cols = @("H","I","J","K","L","M","N")
foreach ($c in $cols) {
    $ws.Range($c + "8").Formula = "=MAX((" + $c + "3/`$C`$3), (" + $c + "4/`$C`$4), (" + $c + "5/`$C`$5), (" + $c + "6/`$C`$6))"
}

# ---------------------------------------------------------------------------
# 2) Row 9 ("Wght Speed"): was (col x/$C$x) summed, now ($C$x/col x) summed
# ---------------------------------------------------------------------------
foreach ($c in $cols) {
    $ws.Range($c + "9").Formula = "=(`$C`$3/" + $c + "3)+(`$C`$4/" + $c + "4)+(`$C`$5/" + $c + "5)+(`$C`$6/" + $c + "6)"
}

# ---------------------------------------------------------------------------
# 3) P:S columns of rows 8 & 9 now simply mirror H:K (instead of being
#    independently recomputed shared formulas)
# ---------------------------------------------------------------------------
$mirrorCols = @{ "P" = "H"; "Q" = "I"; "R" = "J"; "S" = "K" }
foreach ($row in @(8, 9, 10)) {
    foreach ($dst in @("P", "Q", "R", "S")) {
        $src = $mirrorCols[$dst]
        $ws.Range($dst + $row).Formula = "=" + $src + $row
    }
}

# ---------------------------------------------------------------------------
# 4) Insert a new row 11 ("Fairness" = Row9 / Row8) pushing old rows 12-14
#    down to 13-15 (and shifting every chart's series range by one row, fixed
#    up below).
# ---------------------------------------------------------------------------
$ws.Rows("11:11").Insert()

foreach ($c in $cols) {
    $ws.Range($c + "11").Formula = "=" + $c + "9/" + $c + "8"
}

# ---------------------------------------------------------------------------
# 5) Fix up the chart series ranges: every series below old row 11 shifted
#    down by one row when row 11 was inserted.
# ---------------------------------------------------------------------------
$chartFormulas = @(
    "=SERIES(,Sheet1!`$H`$2:`$N`$2,Sheet1!`$H`$13:`$N`$13,1)",
    "=SERIES(,Sheet1!`$H`$2:`$J`$2,Sheet1!`$H`$13:`$J`$13,1)",
    "=SERIES(,Sheet1!`$H`$2:`$J`$2,Sheet1!`$H`$14:`$J`$14,1)",
    "=SERIES(,Sheet1!`$H`$2:`$N`$2,Sheet1!`$H`$14:`$N`$14,1)",
    "=SERIES(,Sheet1!`$H`$2:`$N`$2,Sheet1!`$H`$15:`$N`$15,1)",
    "=SERIES(,Sheet1!`$P`$2:`$S`$2,Sheet1!`$P`$13:`$S`$13,1)",
    "=SERIES(,Sheet1!`$P`$2:`$S`$2,Sheet1!`$P`$14:`$S`$14,1)"
)

$chartCount = $ws.ChartObjects().Count
for ($i = 1; $i -le $chartCount; $i++) {
    $co = $ws.ChartObjects($i)
    $chart = $co.Chart
    $series = $chart.SeriesCollection(1)
    $series.Formula = $chartFormulas[$i - 1]
}

# ---------------------------------------------------------------------------
# 6) New shared-string labels used by the new scratch block below.
# ---------------------------------------------------------------------------
$ws.Range("C63").Value = "IPC alone"
$ws.Range("G63").Value = "IPC shared"

# ---------------------------------------------------------------------------
# 7) New "IPC alone vs IPC shared" scratch calculation block (rows 63-72).
# ---------------------------------------------------------------------------
$baseline = @{ 64 = 1320979; 65 = 13704031; 66 = 6609855; 67 = 1509532 }
$shared   = @{ 64 = 1750297; 65 = 33867573; 66 = 31353825; 67 = 18968058 }

foreach ($row in @(64, 65, 66, 67)) {
    $ws.Range("B" + $row).Value = $baseline[$row]
    $ws.Range("C" + $row).Formula = "=`$F`$64/B" + $row
    $ws.Range("E" + $row).Value = $shared[$row]
    $ws.Range("F" + $row).Value = 5000000
    $ws.Range("G" + $row).Formula = "=F" + $row + "/E" + $row
    $ws.Range("I" + $row).Formula = "=C" + $row + "/G" + $row
}

$ws.Range("E69").Formula = "=E64/B64"
$ws.Range("I69").Formula = "=MAX(I64:I67)"
$ws.Range("E70").Formula = "=E65/B65"
$ws.Range("E71").Formula = "=E66/B66"
$ws.Range("E72").Formula = "=E67/B67"

# ---------------------------------------------------------------------------
# 8) Label row 11 ("Fairness") -- added last so the shared-string table gets
#    appended in the same order as the target workbook (IPC alone, IPC
#    shared, then "Fairness").
# ---------------------------------------------------------------------------
$ws.Range("G11").Value = """Fairness"""
$ws.Range("G11").Font.Bold = $true

# ---------------------------------------------------------------------------
# 9) Restore the saved selection/scroll state.
# ---------------------------------------------------------------------------
$ws.Range("L9").Select()

$wb.Application.Calculate()
